$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value
$ws.Range("B2").Value = 217

# Row 3: name changes to CHUNGA DE LA CRUZ ROSA LILIANA, value 157
$ws.Range("A3").Value = "CHUNGA DE LA CRUZ ROSA LILIANA"
$ws.Range("B3").Value = 157

# Row 4: name changes to ESPINOZA GUZMAN MAYRA LOURDES, value 150
$ws.Range("A4").Value = "ESPINOZA GUZMAN MAYRA LOURDES"
$ws.Range("B4").Value = 150

# Row 5
$ws.Range("B5").Value = 146

# Row 6
$ws.Range("B6").Value = 143

# Row 7
$ws.Range("B7").Value = 116

# Row 8
$ws.Range("B8").Value = 102

# Row 9
$ws.Range("B9").Value = 101

# Row 10
$ws.Range("B10").Value = 96

# Row 11
$ws.Range("B11").Value = 89

# Row 12
$ws.Range("B12").Value = 88

# Row 13
$ws.Range("B13").Value = 75
